$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 636.0855713333334
$ws.Cells.Item(2, 8).Value = 1908.256714
$ws.Cells.Item(2, 9).Value = 0.20580936912678
$ws.Cells.Item(2, 10).Value = 0.20580936912678
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 239.0839323333333
$ws.Cells.Item(2, 14).Value = 717.251797
$ws.Cells.Item(2, 15).Value = 0.4086975387666237
$ws.Cells.Item(2, 16).Value = 0.4086975387666237
$ws.Cells.Item(2, 17).Value = 152077.8396948684
$ws.Cells.Item(2, 18).Value = 1368700.557253815
$ws.Cells.Item(2, 19).Value = 0.08411378261722652
$ws.Cells.Item(2, 20).Value = 0.08411378261722652

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 636.0855713333334
$ws.Cells.Item(3, 8).Value = 1908.256714
$ws.Cells.Item(3, 9).Value = 0.20580936912678
$ws.Cells.Item(3, 10).Value = 0.20580936912678
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 117.0512696666667
$ws.Cells.Item(3, 14).Value = 351.153809
$ws.Cells.Item(3, 15).Value = 0.2000910950200451
$ws.Cells.Item(3, 16).Value = 0.2000910950200451
$ws.Cells.Item(3, 17).Value = 74454.62374121374
$ws.Cells.Item(3, 18).Value = 670091.6136709237
$ws.Cells.Item(3, 19).Value = 0.04118062203396206
$ws.Cells.Item(3, 20).Value = 0.04118062203396206

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 636.0855713333334
$ws.Cells.Item(4, 8).Value = 1908.256714
$ws.Cells.Item(4, 9).Value = 0.20580936912678
$ws.Cells.Item(4, 10).Value = 0.20580936912678
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 171.15883
$ws.Cells.Item(4, 14).Value = 513.47649
$ws.Cells.Item(4, 15).Value = 0.2925842480357353
$ws.Cells.Item(4, 16).Value = 0.2925842480357353
$ws.Cells.Item(4, 17).Value = 108871.6621692949
$ws.Cells.Item(4, 18).Value = 979844.959523654
$ws.Cells.Item(4, 19).Value = 0.06021657950466799
$ws.Cells.Item(4, 20).Value = 0.06021657950466799

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 636.0855713333334
$ws.Cells.Item(5, 8).Value = 1908.256714
$ws.Cells.Item(5, 9).Value = 0.20580936912678
$ws.Cells.Item(5, 10).Value = 0.20580936912678
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 57.695868
$ws.Cells.Item(5, 14).Value = 173.087604
$ws.Cells.Item(5, 15).Value = 0.09862711817759588
$ws.Cells.Item(5, 16).Value = 0.09862711817759588
$ws.Cells.Item(5, 17).Value = 36699.50916035259
$ws.Cells.Item(5, 18).Value = 330295.5824431733
$ws.Cells.Item(5, 19).Value = 0.02029838497092338
$ws.Cells.Item(5, 20).Value = 0.02029838497092338

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1037.896708333333
$ws.Cells.Item(6, 8).Value = 3113.690125
$ws.Cells.Item(6, 9).Value = 0.3358178150670637
$ws.Cells.Item(6, 10).Value = 0.3358178150670637
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 239.0839323333333
$ws.Cells.Item(6, 14).Value = 717.251797
$ws.Cells.Item(6, 15).Value = 0.4086975387666237
$ws.Cells.Item(6, 16).Value = 0.4086975387666237
$ws.Cells.Item(6, 17).Value = 248144.4263841561
$ws.Cells.Item(6, 18).Value = 2233299.837457405
$ws.Cells.Item(6, 19).Value = 0.1372479144918941
$ws.Cells.Item(6, 20).Value = 0.1372479144918941

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1037.896708333333
$ws.Cells.Item(7, 8).Value = 3113.690125
$ws.Cells.Item(7, 9).Value = 0.3358178150670637
$ws.Cells.Item(7, 10).Value = 0.3358178150670637
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 117.0512696666667
$ws.Cells.Item(7, 14).Value = 351.153809
$ws.Cells.Item(7, 15).Value = 0.2000910950200451
$ws.Cells.Item(7, 16).Value = 0.2000910950200451
$ws.Cells.Item(7, 17).Value = 121487.1274932707
$ws.Cells.Item(7, 18).Value = 1093384.147439436
$ws.Cells.Item(7, 19).Value = 0.06719415434400776
$ws.Cells.Item(7, 20).Value = 0.06719415434400777

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1037.896708333333
$ws.Cells.Item(8, 8).Value = 3113.690125
$ws.Cells.Item(8, 9).Value = 0.3358178150670637
$ws.Cells.Item(8, 10).Value = 0.3358178150670637
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 171.15883
$ws.Cells.Item(8, 14).Value = 513.47649
$ws.Cells.Item(8, 15).Value = 0.2925842480357353
$ws.Cells.Item(8, 16).Value = 0.2925842480357353
$ws.Cells.Item(8, 17).Value = 177645.1862591846
$ws.Cells.Item(8, 18).Value = 1598806.676332661
$ws.Cells.Item(8, 19).Value = 0.09825500289840043
$ws.Cells.Item(8, 20).Value = 0.09825500289840046

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1037.896708333333
$ws.Cells.Item(9, 8).Value = 3113.690125
$ws.Cells.Item(9, 9).Value = 0.3358178150670637
$ws.Cells.Item(9, 10).Value = 0.3358178150670637
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 57.695868
$ws.Cells.Item(9, 14).Value = 173.087604
$ws.Cells.Item(9, 15).Value = 0.09862711817759588
$ws.Cells.Item(9, 16).Value = 0.09862711817759588
$ws.Cells.Item(9, 17).Value = 59882.35148163451
$ws.Cells.Item(9, 18).Value = 538941.1633347105
$ws.Cells.Item(9, 19).Value = 0.03312074333276133
$ws.Cells.Item(9, 20).Value = 0.03312074333276133

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 756.5536603333334
$ws.Cells.Item(10, 8).Value = 2269.660981
$ws.Cells.Item(10, 9).Value = 0.2447875546325883
$ws.Cells.Item(10, 10).Value = 0.2447875546325883
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 239.0839323333333
$ws.Cells.Item(10, 14).Value = 717.251797
$ws.Cells.Item(10, 15).Value = 0.4086975387666237
$ws.Cells.Item(10, 16).Value = 0.4086975387666237
$ws.Cells.Item(10, 17).Value = 180879.8241336703
$ws.Cells.Item(10, 18).Value = 1627918.417203033
$ws.Cells.Item(10, 19).Value = 0.1000440710990393
$ws.Cells.Item(10, 20).Value = 0.1000440710990393

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 756.5536603333334
$ws.Cells.Item(11, 8).Value = 2269.660981
$ws.Cells.Item(11, 9).Value = 0.2447875546325883
$ws.Cells.Item(11, 10).Value = 0.2447875546325883
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 117.0512696666667
$ws.Cells.Item(11, 14).Value = 351.153809
$ws.Cells.Item(11, 15).Value = 0.2000910950200451
$ws.Cells.Item(11, 16).Value = 0.2000910950200451
$ws.Cells.Item(11, 17).Value = 88555.56651298075
$ws.Cells.Item(11, 18).Value = 797000.0986168267
$ws.Cells.Item(11, 19).Value = 0.04897980985371371
$ws.Cells.Item(11, 20).Value = 0.04897980985371371

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 756.5536603333334
$ws.Cells.Item(12, 8).Value = 2269.660981
$ws.Cells.Item(12, 9).Value = 0.2447875546325883
$ws.Cells.Item(12, 10).Value = 0.2447875546325883
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 171.15883
$ws.Cells.Item(12, 14).Value = 513.47649
$ws.Cells.Item(12, 15).Value = 0.2925842480357353
$ws.Cells.Item(12, 16).Value = 0.2925842480357353
$ws.Cells.Item(12, 17).Value = 129490.8393348707
$ws.Cells.Item(12, 18).Value = 1165417.554013837
$ws.Cells.Item(12, 19).Value = 0.07162098260068232
$ws.Cells.Item(12, 20).Value = 0.07162098260068232

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 756.5536603333334
$ws.Cells.Item(13, 8).Value = 2269.660981
$ws.Cells.Item(13, 9).Value = 0.2447875546325883
$ws.Cells.Item(13, 10).Value = 0.2447875546325883
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 57.695868
$ws.Cells.Item(13, 14).Value = 173.087604
$ws.Cells.Item(13, 15).Value = 0.09862711817759588
$ws.Cells.Item(13, 16).Value = 0.09862711817759588
$ws.Cells.Item(13, 17).Value = 43650.02012150884
$ws.Cells.Item(13, 18).Value = 392850.1810935795
$ws.Cells.Item(13, 19).Value = 0.024142691079153
$ws.Cells.Item(13, 20).Value = 0.024142691079153

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 660.1181640000001
$ws.Cells.Item(14, 8).Value = 1980.354492
$ws.Cells.Item(14, 9).Value = 0.2135852611735681
$ws.Cells.Item(14, 10).Value = 0.2135852611735681
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 239.0839323333333
$ws.Cells.Item(14, 14).Value = 717.251797
$ws.Cells.Item(14, 15).Value = 0.4086975387666237
$ws.Cells.Item(14, 16).Value = 0.4086975387666237
$ws.Cells.Item(14, 17).Value = 157823.6464537803
$ws.Cells.Item(14, 18).Value = 1420412.818084022
$ws.Cells.Item(14, 19).Value = 0.08729177055846378
$ws.Cells.Item(14, 20).Value = 0.08729177055846379

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 660.1181640000001
$ws.Cells.Item(15, 8).Value = 1980.354492
$ws.Cells.Item(15, 9).Value = 0.2135852611735681
$ws.Cells.Item(15, 10).Value = 0.2135852611735681
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 117.0512696666667
$ws.Cells.Item(15, 14).Value = 351.153809
$ws.Cells.Item(15, 15).Value = 0.2000910950200451
$ws.Cells.Item(15, 16).Value = 0.2000910950200451
$ws.Cells.Item(15, 17).Value = 77267.6692262289
$ws.Cells.Item(15, 18).Value = 695409.0230360602
$ws.Cells.Item(15, 19).Value = 0.04273650878836156
$ws.Cells.Item(15, 20).Value = 0.04273650878836156

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 660.1181640000001
$ws.Cells.Item(16, 8).Value = 1980.354492
$ws.Cells.Item(16, 9).Value = 0.2135852611735681
$ws.Cells.Item(16, 10).Value = 0.2135852611735681
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 171.15883
$ws.Cells.Item(16, 14).Value = 513.47649
$ws.Cells.Item(16, 15).Value = 0.2925842480357353
$ws.Cells.Item(16, 16).Value = 0.2925842480357353
$ws.Cells.Item(16, 17).Value = 112985.0526119881
$ws.Cells.Item(16, 18).Value = 1016865.473507893
$ws.Cells.Item(16, 19).Value = 0.06249168303198455
$ws.Cells.Item(16, 20).Value = 0.06249168303198455

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 660.1181640000001
$ws.Cells.Item(17, 8).Value = 1980.354492
$ws.Cells.Item(17, 9).Value = 0.2135852611735681
$ws.Cells.Item(17, 10).Value = 0.2135852611735681
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 57.695868
$ws.Cells.Item(17, 14).Value = 173.087604
$ws.Cells.Item(17, 15).Value = 0.09862711817759588
$ws.Cells.Item(17, 16).Value = 0.09862711817759588
$ws.Cells.Item(17, 17).Value = 38086.09045454636
$ws.Cells.Item(17, 18).Value = 342774.8140909172
$ws.Cells.Item(17, 19).Value = 0.02106529879475818
$ws.Cells.Item(17, 20).Value = 0.02106529879475818
